# edit.ps1 — apply "workflow vs validation added" commit
#
# 1) Insert a new first paragraph (bold, yellow-highlighted) in front of the
#    existing "Q: Can a workflow rule be triggered..." paragraph, explaining
#    that a workflow rule does not respect a Validation Rule.
# 2) Insert two empty (bold-marked) paragraphs right before the existing
#    "Q: What are assignment Rules?" heading paragraph.

$d = $word.ActiveDocument

# --- 1. New opening paragraph -------------------------------------------------
$firstParaRange = $d.Paragraphs(1).Range
$firstParaRange.InsertParagraphBefore()

$newFirstPara = $d.Paragraphs(1).Range
$newFirstPara.Text = "A workflow rule does not does not respect Validation Rule. So, if there is a situation where workflow and validation rule are conflicting, workflow rule will win because of the order of execution prescribed by Salesforce."
$newFirstPara.Font.Bold = $true
$newFirstPara.Font.BoldBi = $true
$newFirstPara.HighlightColorIndex = 7

# --- 2. Two blank paragraphs before "Q: What are assignment Rules?" ---------
$target = $d.Content
$target.Find.Execute("Q: What are assignment Rules?") | Out-Null
$targetPara = $target.Paragraphs(1).Range
$targetPara.InsertParagraphBefore()
$targetPara.InsertParagraphBefore()
